$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.375.81"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "3.328.41"
$ws.Range("E3").Value = "  -4.29%  "
$c = $ws.Range("D5")
$c.Value = "'574.64"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.32%  "
$c = $ws.Range("D6")
$c.Value = "'177.13"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.76%  "
$ws.Range("E7").Value = "  +3.32%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "3.327.59"
$ws.Range("E9").Value = "  -4.27%  "
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "3.908.99"
$ws.Range("E13").Value = "  -4.20%  "
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("E15").Value = "  -4.56%  "
$ws.Range("D16").Value = "65.558.64"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "3.324.28"
$ws.Range("E18").Value = "  -4.46%  "
$c = $ws.Range("D19")
$c.Value = "'5.72"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -3.46%  "
$c = $ws.Range("D20")
$c.Value = "'13.37"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'362.30"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.01%  "
$c = $ws.Range("D22")
$c.Value = "'7.42"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -4.32%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D23")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D24")
$c.Value = "'71.23"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("E25").Value = "  -3.21%  "
$c = $ws.Range("D26")
$c.Value = "'0.0000121"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.93%  "
$c = $ws.Range("D27")
$c.Value = "'9.59"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.02%  "
$c = $ws.Range("D28")
$c.Value = "'0.178"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -1.55%  "
$c = $ws.Range("D31")
$c.Value = "'5.63"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -5.07%  "
$ws.Range("E34").Value = "  -4.16%  "
$ws.Range("E35").Value = "  -5.29%  "
$c = $ws.Range("D36")
$c.Value = "'1.50"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.42%  "
$c = $ws.Range("D37")
$c.Value = "'159.87"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "
$c = $ws.Range("D38")
$c.Value = "'0.845"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.99%  "
$c = $ws.Range("D39")
$c.Value = "'27.14"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -7.16%  "
$ws.Range("E40").Value = "  -1.18%  "
$c = $ws.Range("D41")
$c.Value = "'2.53"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.12%  "
$ws.Range("D42").Value = "2.686.24"
$ws.Range("E42").Value = "  -4.93%  "
$c = $ws.Range("D43")
$c.Value = "'6.21"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -3.37%  "
$c = $ws.Range("D44")
$c.Value = "'4.28"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.87%  "
$c = $ws.Range("D45")
$c.Value = "'0.0665"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.56%  "
$c = $ws.Range("D46")
$c.Value = "'39.65"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.68%  "
$c = $ws.Range("D47")
$c.Value = "'333.76"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +7.63%  "
$c = $ws.Range("D48")
$c.Value = "'24.28"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("E49").Value = "  -3.48%  "
$ws.Range("E50").Value = "  +2.36%  "
$ws.Range("E51").Value = "  -1.37%  "
